# 451: add ct_guarantee and ct_guarantee_ff fields
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# Populate cells in the same order the values were first authored, so that
# new shared-string table entries are appended in the expected sequence.
$ws.Range("F4").Value = "1;2;3;4;5;6;7;8"
$ws.Range("F1").Value = "ct_guarantee"
$ws.Range("F2").Value = "1;2;3"
$ws.Range("F3").Value = "1;2;3"
$ws.Range("F5").Value = "10000;1200;1;2"
$ws.Range("G1").Value = "ct_guarantee_ff"
$ws.Range("G2").Value = "abc;def"
$ws.Range("G3").Value = "abc;def"
$ws.Range("G4").Value = "abc;def"
$ws.Range("G5").Value = "abc;def"
$ws.Range("G6").Value = "abc;def"
$ws.Range("G9").Value = "abc;def"
$ws.Range("G10").Value = "abc;def"
$ws.Range("G11").Value = "abc;def"
$ws.Range("F7").Value = "1;2;3;4;5;6"
$ws.Range("F10").Value = "977;1"

# Numeric values (do not touch the shared-string table)
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = 977
$ws.Range("F9").Value = 977
$ws.Range("F11").Value = 999

# Column widths for new columns
$ws.Range("F1").EntireColumn.ColumnWidth = 19.83203125
$ws.Range("G1").EntireColumn.ColumnWidth = 14.6640625

# Row height for header row (wraps to two lines now)
$ws.Rows.Item(1).RowHeight = 34

# View adjustments
$ws.Range("E14").Select()
$excel.ActiveWindow.ScrollColumn = 3
